# This script updates the "Pais" (countries) COVID-19 stats sheet to match a newer
# data snapshot (timestamp bumped from 14:29 to 15:29) that also introduced three
# countries that were missing before: Serbia, Republica de Macedonia and Mauricio.
# Each new country is inserted directly above its former alphabetical neighbour in the
# (descending-by-total-cases-sorted) table, which pushes the rows that used to sit at
# that position down by one; three already-listed countries (Estados Unidos, Reino
# Unido, Austria) simply get refreshed figures. All other rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 15:29"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 104996
$ws.Range("C4").Value = 870
$ws.Range("D4").Value = 2537
$ws.Range("E4").Value = 100742
$ws.Range("F4").Value = 2494
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 1717

# Row 11: Reino Unido
$ws.Range("A11").Value = "Reino Unido"
$ws.Range("B11").Value = 17089
$ws.Range("C11").Value = 2546
$ws.Range("D11").Value = 135
$ws.Range("E11").Value = 15935
$ws.Range("F11").Value = 163
$ws.Range("G11").Value = 260
$ws.Range("H11").Value = 1019

# Row 16: Austria
$ws.Range("A16").Value = "Austria"
$ws.Range("B16").Value = 8030
$ws.Range("C16").Value = 333
$ws.Range("D16").Value = 225
$ws.Range("E16").Value = 7737
$ws.Range("F16").Value = 128
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 68

# Row 52: Serbia
$ws.Range("A52").Value = "Serbia"
$ws.Range("B52").Value = 659
$ws.Range("C52").Value = 131
$ws.Range("D52").Value = 42
$ws.Range("E52").Value = 607
$ws.Range("F52").Value = 25
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 10

# Row 53: Estonia
$ws.Range("A53").Value = "Estonia"
$ws.Range("B53").Value = 645
$ws.Range("C53").Value = 70
$ws.Range("D53").Value = 20
$ws.Range("E53").Value = 624
$ws.Range("F53").Value = 10
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 1

# Row 54: Peru
$ws.Range("A54").Value = "Peru"
$ws.Range("B54").Value = 635
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 16
$ws.Range("E54").Value = 608
$ws.Range("F54").Value = 21
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 11

# Row 55: Croacia
$ws.Range("A55").Value = "Croacia"
$ws.Range("B55").Value = 635
$ws.Range("C55").Value = 49
$ws.Range("D55").Value = 45
$ws.Range("E55").Value = 586
$ws.Range("F55").Value = 14
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 4

# Row 56: Republica Dominicana
$ws.Range("A56").Value = "Republica Dominicana"
$ws.Range("B56").Value = 581
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 3
$ws.Range("E56").Value = 558
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 20

# Row 57: Catar
$ws.Range("A57").Value = "Catar"
$ws.Range("B57").Value = 562
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 43
$ws.Range("E57").Value = 519
$ws.Range("F57").Value = 6
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0

# Row 58: Hong Kong
$ws.Range("A58").Value = "Hong Kong"
$ws.Range("B58").Value = 560
$ws.Range("C58").Value = 41
$ws.Range("D58").Value = 112
$ws.Range("E58").Value = 444
$ws.Range("F58").Value = 5
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 4

# Row 59: Colombia
$ws.Range("A59").Value = "Colombia"
$ws.Range("B59").Value = 539
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 10
$ws.Range("E59").Value = 523
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 6

# Row 60: Egipto
$ws.Range("A60").Value = "Egipto"
$ws.Range("B60").Value = 536
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 116
$ws.Range("E60").Value = 390
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 30

# Row 80: Republica de Macedonia
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 241
$ws.Range("C80").Value = 22
$ws.Range("D80").Value = 3
$ws.Range("E80").Value = 234
$ws.Range("F80").Value = 1
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 4

# Row 81: Jordania
$ws.Range("A81").Value = "Jordania"
$ws.Range("B81").Value = 235
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 18
$ws.Range("E81").Value = 216
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 1

# Row 82: Kuwait
$ws.Range("A82").Value = "Kuwait"
$ws.Range("B82").Value = 235
$ws.Range("C82").Value = 10
$ws.Range("D82").Value = 64
$ws.Range("E82").Value = 171
$ws.Range("F82").Value = 11
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 0

# Row 83: Tunez
$ws.Range("A83").Value = "Tunez"
$ws.Range("B83").Value = 227
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 2
$ws.Range("E83").Value = 218
$ws.Range("F83").Value = 10
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 7

# Row 84: San Marino
$ws.Range("A84").Value = "San Marino"
$ws.Range("B84").Value = 223
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 6
$ws.Range("E84").Value = 196
$ws.Range("F84").Value = 15
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 21

# Row 103: Mauricio
$ws.Range("A103").Value = "Mauricio"
$ws.Range("B103").Value = 102
$ws.Range("C103").Value = 8
$ws.Range("D103").Value = 0
$ws.Range("E103").Value = 100
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 2

# Row 104: Costa de Marfil
$ws.Range("A104").Value = "Costa de Marfil"
$ws.Range("B104").Value = 101
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 3
$ws.Range("E104").Value = 98
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

# Row 105: Camboya
$ws.Range("A105").Value = "Camboya"
$ws.Range("B105").Value = 99
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 13
$ws.Range("E105").Value = 86
$ws.Range("F105").Value = 1
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

# Row 106: Estado de Palestina
$ws.Range("A106").Value = "Estado de Palestina"
$ws.Range("B106").Value = 97
$ws.Range("C106").Value = 6
$ws.Range("D106").Value = 18
$ws.Range("E106").Value = 78
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 1

# Row 107: Guadalupe
$ws.Range("A107").Value = "Guadalupe"
$ws.Range("B107").Value = 96
$ws.Range("C107").Value = 23
$ws.Range("D107").Value = 17
$ws.Range("E107").Value = 77
$ws.Range("F107").Value = 4
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 2

# Row 108: Honduras
$ws.Range("A108").Value = "Honduras"
$ws.Range("B108").Value = 95
$ws.Range("C108").Value = 27
$ws.Range("D108").Value = 3
$ws.Range("E108").Value = 91
$ws.Range("F108").Value = 4
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 1
